$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-12-05 Tuesday" "2023-12-06 Wednesday"

Replace-Text "17×51=867" "50×40=2000"
Replace-Text "74×49=3626" "43×19=817"
Replace-Text "68×25=1700" "73×15=1095"
Replace-Text "55×45=2475" "96×93=8928"
Replace-Text "63×56=3528" "91×90=8190"

Replace-Text "31×56=1736" "98×67=6566"
Replace-Text "13×48=624" "36×36=1296"
Replace-Text "88×79=6952" "87×48=4176"
Replace-Text "98×17=1666" "60×79=4740"
Replace-Text "55×40=2200" "16×61=976"

Replace-Text "47×83=3901" "68×13=884"
Replace-Text "75×56=4200" "29×69=2001"
Replace-Text "92×60=5520" "50×94=4700"
Replace-Text "44×45=1980" "79×77=6083"
Replace-Text "97×25=2425" "59×75=4425"

Replace-Text "68×51=3468" "67×51=3417"
Replace-Text "45×49=2205" "94×91=8554"
Replace-Text "37×98=3626" "93×15=1395"
Replace-Text "44×79=3476" "70×79=5530"
Replace-Text "89×39=3471" "93×64=5952"

Replace-Text "70×73=5110" "30×27=810"
Replace-Text "75×24=1800" "69×66=4554"
Replace-Text "39×44=1716" "43×71=3053"
Replace-Text "21×12=252" "65×18=1170"
Replace-Text "99×76=7524" "64×75=4800"
